# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) for a batch of leve rows across all profession
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Mirrors the upstream
# runner's output — some rows also lose a stale cell (ClearContents) where
# the refreshed data no longer carries that figure.

$wb = $excel.ActiveWorkbook

$colIndex = @{ "H" = 8; "I" = 9; "J" = 10; "K" = 11; "L" = 12; "M" = 13; "N" = 14 }

function Set-RowValues {
    param(
        $ws,
        [int]$row,
        [hashtable]$values
    )
    foreach ($colLetter in $values.Keys) {
        $col = $colIndex[$colLetter]
        $val = $values[$colLetter]
        if ($null -eq $val) {
            $ws.Cells.Item($row, $col).ClearContents()
        } else {
            $ws.Cells.Item($row, $col).Value = $val
        }
    }
}

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")

Set-RowValues $ws 43 @{
    H = 1724.4; I = 1823.3334; J = 1682; K = 1823.3334; L = 1682;
    M = -1754.3334; N = -1820
}
Set-RowValues $ws 129 @{
    H = 995.5909; I = 440.3; J = 1094.75; K = 1320.9; L = 3284.25;
    M = 3679.1; N = -13284.25
}
Set-RowValues $ws 135 @{
    H = 13889856; I = 16667226; K = 150005034; M = -150002499
}
Set-RowValues $ws 137 @{
    H = 1603925.4; I = 2977351.8; J = 1594.75; K = 8932055.399999999;
    L = 4784.25; M = -8929505.399999999; N = -9884.25
}
Set-RowValues $ws 138 @{
    H = 2979.69; I = 1571.6; J = 3228.1765; K = 4714.799999999999;
    L = 9684.529500000001; M = 425.2000000000007; N = -19964.5295
}

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")

Set-RowValues $ws 2 @{
    H = 1660; I = 1480; J = 1750; K = 1480; L = 1750; M = -1367; N = -1976
}
Set-RowValues $ws 32 @{
    H = 6863175.5; I = 8208840; K = 8208840; M = -8208553
}
Set-RowValues $ws 45 @{
    H = 2258.9443; I = 1844.7273; J = 2909.8572; K = 1844.7273;
    L = 2909.8572; M = -1467.7273; N = -3663.8572
}
Set-RowValues $ws 46 @{
    H = 2976; J = 2976; L = 2976; N = -3614
}
Set-RowValues $ws 63 @{
    H = 4708.6; I = 4300; J = 4737.7856; K = 4300; L = 4737.7856;
    M = -3614; N = -6109.7856
}
Set-RowValues $ws 66 @{
    H = 4708.6; I = 4300; J = 4737.7856; K = 21500; L = 23688.928;
    M = -18068; N = -30552.928
}
Set-RowValues $ws 102 @{
    H = 2316; J = 2666.6667; L = 2666.6667; N = -5910.6667
}
Set-RowValues $ws 116 @{
    H = 1660; I = 1480; J = 1750; K = 1480; L = 1750; M = 814; N = -6338
}
Set-RowValues $ws 132 @{
    H = 3092.672; I = 3272; J = 2851.2693; K = 9816; L = 8553.8079;
    M = -7286; N = -13613.8079
}

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")

Set-RowValues $ws 3 @{
    H = 1660; I = 1480; J = 1750; K = 1480; L = 1750; M = -1366; N = -1978
}
Set-RowValues $ws 86 @{
    H = 2194.5715; I = 2151; J = 2303.5; K = 2151; L = 2303.5;
    M = -1028; N = -4549.5
}
Set-RowValues $ws 89 @{
    H = 2194.5715; I = 2151; J = 2303.5; K = 10755; L = 11517.5;
    M = -5139; N = -22749.5
}
Set-RowValues $ws 134 @{
    H = 3221.0435; I = 3002.6155; J = 3505; K = 9007.8465; L = 10515;
    M = -6472.8465; N = -15585
}

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")

Set-RowValues $ws 99 @{
    H = 2436.6592; I = 2255.3333; J = 2483.2856; K = 2255.3333;
    L = 2483.2856; M = -757.3332999999998; N = -5479.2856
}
Set-RowValues $ws 122 @{
    H = 1569.1714; I = 1158.2778; J = 2004.2354; K = 3474.8334;
    L = 6012.706200000001; M = -1024.8334; N = -10912.7062
}
Set-RowValues $ws 126 @{
    H = 2436.6592; I = 2255.3333; J = 2483.2856; K = 6765.999899999999;
    L = 7449.8568; M = -4295.999899999999; N = -12389.8568
}

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")

Set-RowValues $ws 98 @{
    H = 143411; I = 600; J = 200535.4; K = 1800; L = 601606.2;
    M = -302; N = -604602.2
}
Set-RowValues $ws 104 @{
    H = 3818.4285; I = 0; J = 3818.4285; K = 0; L = 11455.2855;
    M = $null; N = -16697.2855
}
Set-RowValues $ws 113 @{
    H = 573.9231; I = 592.3077; J = 555.53845; K = 1776.9231;
    L = 1666.61535; M = 393.0769; N = -6006.61535
}

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")

Set-RowValues $ws 122 @{
    H = 12548.9; I = 18949.834; J = 2947.5; K = 56849.50199999999;
    L = 8842.5; M = -54399.50199999999; N = -13742.5
}
Set-RowValues $ws 123 @{
    H = 12353.708; J = 14224.45; L = 14224.45; N = -19124.45
}
Set-RowValues $ws 126 @{
    H = 3350; I = 3350; J = 0; K = 10050; L = 0; M = -7580; N = $null
}

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")

Set-RowValues $ws 55 @{
    H = 707.5; I = 470; J = 945; K = 470; L = 945; M = -297; N = -1291
}
Set-RowValues $ws 122 @{
    H = 5591.95; I = 2600; K = 7800; M = -5350
}

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")

Set-RowValues $ws 122 @{
    H = 2744.889; I = 2775.5; K = 8326.5; M = -5876.5
}
Set-RowValues $ws 132 @{
    H = 3706295.2; I = 2466.1292; J = 11907631; K = 7398.3876;
    L = 35722893; M = -4868.3876; N = -35727953
}

Write-Output "Leve price refresh applied."
